$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 108, pushing the existing rows 108-128 down to 110-130.
$ws.Rows("108:109").Insert()

# Row 108 - new "Early Majestic" / Primera record
$ws.Cells.Item(108, 1).Value = 10
$ws.Cells.Item(108, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(108, 3).Value = "La Araucanía"
$ws.Cells.Item(108, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 4 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Cells.Item(108, 5).Value = 9
$ws.Cells.Item(108, 6).Value = "Fruta"
$ws.Cells.Item(108, 7).Value = 100103
$ws.Cells.Item(108, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(108, 9).Value = 100103004
$ws.Cells.Item(108, 10).Value = "Durazno"
$ws.Cells.Item(108, 11).Value = "Early Majestic"
$ws.Cells.Item(108, 12).Value = "Primera"
$ws.Cells.Item(108, 13).Value = 80
$ws.Cells.Item(108, 14).Value = 32000
$ws.Cells.Item(108, 15).Value = 32000
$ws.Cells.Item(108, 16).Value = 32000
$ws.Cells.Item(108, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(108, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(108, 19).Value = 1778
$ws.Cells.Item(108, 20).Value = 18

# Row 109 - new "Early Majestic" / Segunda record
$ws.Cells.Item(109, 1).Value = 10
$ws.Cells.Item(109, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(109, 3).Value = "La Araucanía"
$ws.Cells.Item(109, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 4 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Cells.Item(109, 5).Value = 9
$ws.Cells.Item(109, 6).Value = "Fruta"
$ws.Cells.Item(109, 7).Value = 100103
$ws.Cells.Item(109, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(109, 9).Value = 100103004
$ws.Cells.Item(109, 10).Value = "Durazno"
$ws.Cells.Item(109, 11).Value = "Early Majestic"
$ws.Cells.Item(109, 12).Value = "Segunda"
$ws.Cells.Item(109, 13).Value = 200
$ws.Cells.Item(109, 14).Value = 10000
$ws.Cells.Item(109, 15).Value = 10000
$ws.Cells.Item(109, 16).Value = 10000
$ws.Cells.Item(109, 17).Value = "`$/bandeja 8 kilos granel"
$ws.Cells.Item(109, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(109, 19).Value = 1250
$ws.Cells.Item(109, 20).Value = 8
